$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "24.647.31"
Set-TextCell "E2" "  +0.22%  "
Set-TextCell "D3" "1.696.83"
Set-TextCell "E3" "  +0.03%  "
Set-TextCell "D5" "315.50"
Set-TextCell "E5" "  -0.34%  "
Set-TextCell "E6" "  +0.03%  "
Set-TextCell "D7" "0.3922"
Set-TextCell "E7" "  -0.46%  "
Set-TextCell "D8" "0.4046"
Set-TextCell "E8" "  +0.49%  "
Set-TextCell "D9" "1.524"
Set-TextCell "E9" "  -0.07%  "
Set-TextCell "D10" "1.003"
Set-TextCell "E10" "  +0.19%  "
Set-TextCell "D11" "53.10"
Set-TextCell "E11" "  -1.40%  "
Set-TextCell "D12" "0.08838"
Set-TextCell "E12" "  +0.64%  "
Set-TextCell "D13" "7.431"
Set-TextCell "E13" "  +2.71%  "
Set-TextCell "D14" "23.59"
Set-TextCell "E14" "  +1.43%  "
Set-TextCell "D15" "8.129"
Set-TextCell "E15" "  +6.89%  "
Set-TextCell "E16" "  -0.33%  "
Set-TextCell "D17" "1.697.90"
Set-TextCell "E17" "  -0.03%  "
Set-TextCell "D18" "99.27"
Set-TextCell "E18" "  -1.31%  "
Set-TextCell "D19" "0.07017"
Set-TextCell "E19" "  -0.60%  "
Set-TextCell "D20" "19.72"
Set-TextCell "E20" "  +0.14%  "
Set-TextCell "D21" "7.057"
Set-TextCell "E21" "  +2.78%  "
Set-TextCell "D22" "1.004"
Set-TextCell "E22" "  +0.36%  "
Set-TextCell "D23" "14.68"
Set-TextCell "E23" "  +4.31%  "
Set-TextCell "D24" "24.612.96"
Set-TextCell "E24" "  +0.14%  "
Set-TextCell "D25" "3.137"
Set-TextCell "E25" "  +3.41%  "
Set-TextCell "E26" "  +1.58%  "
Set-TextCell "D27" "22.62"
Set-TextCell "E27" "  +0.96%  "
Set-TextCell "D28" "163.50"
Set-TextCell "E28" "  +2.08%  "
Set-TextCell "D29" "8.845"
Set-TextCell "E29" "  +18.28%  "
Set-TextCell "D30" "135.50"
Set-TextCell "E30" "  +0.60%  "
Set-TextCell "D31" "5.139"
Set-TextCell "E31" "  -1.63%  "
Set-TextCell "D32" "0.08994"
Set-TextCell "E32" "  +5.46%  "
Set-TextCell "D33" "7.593"
Set-TextCell "E33" "  +4.15%  "
Set-TextCell "D34" "1.069"
Set-TextCell "E34" "  -3.72%  "
Set-TextCell "D35" "1.963"
Set-TextCell "E35" "  +0.48%  "
Set-TextCell "D36" "11.03"
Set-TextCell "E36" "  -3.16%  "
Set-TextCell "D37" "0.2752"
Set-TextCell "E37" "  +0.47%  "
Set-TextCell "D38" "0.02912"
Set-TextCell "E38" "  +5.37%  "
Set-TextCell "D39" "14.41"
Set-TextCell "E39" "  -0.86%  "
Set-TextCell "D40" "0.09154"
Set-TextCell "E40" "  +0.95%  "
Set-TextCell "D41" "1.454"
Set-TextCell "E41" "  -0.73%  "
Set-TextCell "D42" "0.7657"
Set-TextCell "E42" "  -0.81%  "
Set-TextCell "D43" "15.94"
Set-TextCell "E43" "  +2.57%  "
Set-TextCell "D44" "0.7174"
Set-TextCell "E44" "  -0.40%  "
Set-TextCell "D45" "2.584"
Set-TextCell "E45" "  +1.42%  "
Set-TextCell "E46" "  -0.38%  "
Set-TextCell "E47" "  +0.04%  "
Set-TextCell "D48" "1.334"
Set-TextCell "E48" "  -2.18%  "
Set-TextCell "D49" "139.84"
Set-TextCell "E49" "  -1.01%  "
Set-TextCell "D50" "0.07969"
Set-TextCell "E50" "  -0.69%  "
Set-TextCell "D51" "90.38"
Set-TextCell "E51" "  +2.10%  "
